$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.960.88"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +0.28%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.918.95"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +0.61%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'592.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +1.10%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'145.65"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +0.19%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'  +1.04%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'6.87"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +2.57%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  +0.36%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.439"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  -1.49%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.0000226"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +1.16%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'33.57"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +0.54%  "
$ws.Range("E13").ClearFormats()
$ws.Range("E14").Value = "'  +0.01%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'3.399.09"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +0.55%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'60.897.35"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +0.40%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'6.69"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -0.60%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'2.916.59"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +0.52%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'431.21"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +1.71%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'13.37"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -1.06%  "
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'  +1.74%  "
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = "'  +0.04%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'81.45"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +1.98%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'10.94"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -0.48%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'2.21"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +0.21%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'11.95"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +1.30%  "
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'  +0.06%  "
$ws.Range("E27").ClearFormats()
$ws.Range("E28").Value = "'  +6.42%  "
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = "'  -0.09%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'  +0.09%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'7.05"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -1.50%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'26.47"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +0.89%  "
$ws.Range("E32").ClearFormats()
$ws.Range("E33").Value = "'  +1.24%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'0.0₃0855"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +2.52%  "
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'  +1.20%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'5.63"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +0.78%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'3.03"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +2.58%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'1.99"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -1.05%  "
$ws.Range("E38").ClearFormats()
$ws.Range("E39").Value = "'  -1.78%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'8.55"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -1.19%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.287"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -0.50%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'39.99"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -2.70%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'375.02"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +0.81%  "
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = "'  -0.20%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'2.703.60"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +1.63%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'131.59"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -1.01%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D48").Value = "'23.88"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -4.97%  "
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "'  +0.38%  "
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'  -2.94%  "
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'  +1.98%  "
$ws.Range("E51").ClearFormats()
